# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet, matching the freshly re-scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 767
$ws1.Range("F5").Value = 37
$ws1.Range("F7").Value = 3667
$ws1.Range("F9").Value = 4273
$ws1.Range("F11").Value = 1065

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 767
$ws4.Range("F5").Value = 37
$ws4.Range("F8").Value = 3667
$ws4.Range("F10").Value = 4273
$ws4.Range("F12").Value = 1065
